$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the long "Assessment text" string (K2): add spaces before "(PMID" occurrences ---
$oldText = $ws.Range("K2").Value2
$newText = $oldText.Replace("gain-of-function(PMID", "gain-of-function (PMID").Replace("evidence(PMID:34850743", "evidence (PMID:34850743")
$ws.Range("K2").Value = $newText

# --- Row 3: DYSF variant ---
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = "DYSF"
$ws.Cells.Item(3, 3).Value = "NM_003494.4"
$ws.Cells.Item(3, 4).Value = "c.237-159_342+1237del"
$ws.Cells.Item(3, 5).Value = "p.Phe80ProfsTer36"

$ws.Cells.Item(2, 8).Copy() | Out-Null
$ws.Cells.Item(3, 8).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(3, 8).Value = "Antisense Oligonucleotide"

$ws.Cells.Item(3, 9).Value = "Exon skipping"
$ws.Cells.Item(3, 10).Value = "Not eligible"
$ws.Rows.Item(3).RowHeight = 15.75

# --- Row 4: SCN1A variant ---
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = "SCN1A"

$ws.Cells.Item(4, 3).Value = "`nNM_001165963.4c.4465C>A (p.Gln1489Lys)"
$ws.Cells.Item(4, 3).WrapText = $true

$ws.Cells.Item(4, 4).Value = "c.4465C>A "
$ws.Cells.Item(4, 5).Value = "p.Gln1489Lys"

$ws.Cells.Item(2, 8).Copy() | Out-Null
$ws.Cells.Item(4, 8).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(4, 8).Value = "Antisense Oligonucleotide"

$ws.Cells.Item(4, 10).Value = "Unable to assess"
$ws.Cells.Item(4, 11).Value = "Pathomechanism unknown"
$ws.Rows.Item(4).RowHeight = 15.75

# --- Column width tweaks ---
$ws.Columns.Item(4).ColumnWidth = 20.666666666666664
$ws.Columns.Item(8).ColumnWidth = 23.666666666666664

# --- Selection / view state ---
$ws.Range("K12").Select() | Out-Null
